$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.20133293946381
$ws.Range("C2").Value = 9.666290596202549
$ws.Range("E2").Value = 23.04931394281904
$ws.Range("F2").Value = 36.29890899816973
$ws.Range("G2").Value = 19.19231553601168
$ws.Range("H2").Value = 11.96103799314165
$ws.Range("J2").Value = 7.353762464545818
$ws.Range("M2").Value = 19.25285945647081
$ws.Range("O2").Value = 16.9125619321515
$ws.Range("B3").Value = 11.51876793840451
$ws.Range("C3").Value = 9.304127779162064
$ws.Range("E3").Value = 23.01525262893884
$ws.Range("F3").Value = 36.3105425107217
$ws.Range("G3").Value = 19.33913346232814
$ws.Range("H3").Value = 12.02684006893154
$ws.Range("J3").Value = 7.369982786320945
$ws.Range("M3").Value = 18.96946233524631
$ws.Range("O3").Value = 17.02959569345609
$ws.Range("B4").Value = 11.07768141459265
$ws.Range("C4").Value = 9.074005827163695
$ws.Range("E4").Value = 22.99889935178487
$ws.Range("F4").Value = 36.32878974363857
$ws.Range("G4").Value = 19.44047000871002
$ws.Range("H4").Value = 12.06979353088736
$ws.Range("J4").Value = 7.380520774516295
$ws.Range("M4").Value = 18.79526488931599
$ws.Range("O4").Value = 17.10677617064804
$ws.Range("B5").Value = 10.89252774102524
$ws.Range("C5").Value = 8.978394932657146
$ws.Range("E5").Value = 22.99338649497118
$ws.Range("F5").Value = 36.33901298429196
$ws.Range("G5").Value = 19.48454828870362
$ws.Range("H5").Value = 12.08793877666873
$ws.Range("J5").Value = 7.384960979582359
$ws.Range("M5").Value = 18.72430565288426
$ws.Range("O5").Value = 17.13956161642744
$ws.Range("B6").Value = 10.86146061340459
$ws.Range("C6").Value = 8.962411651826404
$ws.Range("E6").Value = 22.99254073473151
$ws.Range("F6").Value = 36.34087871258713
$ws.Range("G6").Value = 19.49203461542497
$ws.Range("H6").Value = 12.09099051708292
$ws.Range("J6").Value = 7.385707094932393
$ws.Range("M6").Value = 18.71252681308554
$ws.Range("O6").Value = 17.14508604746823
$ws.Range("B7").Value = 11.07520608103089
$ws.Range("C7").Value = 9.07272365051266
$ws.Range("E7").Value = 22.9988203369243
$ws.Range("F7").Value = 36.32891633996458
$ws.Range("G7").Value = 19.44105323755358
$ws.Range("H7").Value = 12.0700356472063
$ws.Range("J7").Value = 7.380580065466583
$ws.Range("M7").Value = 18.79430769614994
$ws.Range("O7").Value = 17.10721293186901
$ws.Range("B8").Value = 11.97061365969438
$ws.Range("C8").Value = 9.543090903133573
$ws.Range("E8").Value = 23.0366259916059
$ws.Range("F8").Value = 36.30061319625707
$ws.Range("G8").Value = 19.24059747273348
$ws.Range("H8").Value = 11.98319711910449
$ws.Range("J8").Value = 7.359235365052057
$ws.Range("M8").Value = 19.15522876177392
$ws.Range("O8").Value = 16.95180813305623
$ws.Range("B9").Value = 13.54787616076933
$ws.Range("C9").Value = 10.39955985893399
$ws.Range("E9").Value = 23.14671774528086
$ws.Range("F9").Value = 36.33334922288818
$ws.Range("G9").Value = 18.93764456103746
$ws.Range("H9").Value = 11.83314760302594
$ws.Range("J9").Value = 7.321952332463384
$ws.Range("M9").Value = 19.85820908590002
$ws.Range("O9").Value = 16.68947706002279
$ws.Range("B10").Value = 14.59358416931043
$ws.Range("C10").Value = 10.98355835552616
$ws.Range("E10").Value = 23.24916201341
$ws.Range("F10").Value = 36.41124394612019
$ws.Range("G10").Value = 18.77181815569832
$ws.Range("H10").Value = 11.73524797905469
$ws.Range("J10").Value = 7.297324595127806
$ws.Range("M10").Value = 20.36748302709461
$ws.Range("O10").Value = 16.52287321452948
$ws.Range("B11").Value = 15.04416157535256
$ws.Range("C11").Value = 11.23855511573931
$ws.Range("E11").Value = 23.30035399149222
$ws.Range("F11").Value = 36.4583382507132
$ws.Range("G11").Value = 18.70906997730876
$ws.Range("H11").Value = 11.69339214604798
$ws.Range("J11").Value = 7.286715875278932
$ws.Range("M11").Value = 20.59672790069299
$ws.Range("O11").Value = 16.45281556876588
$ws.Range("B12").Value = 15.2111372159645
$ws.Range("C12").Value = 11.3335231562166
$ws.Range("E12").Value = 23.32038965623305
$ws.Range("F12").Value = 36.47784173285903
$ws.Range("G12").Value = 18.68716184685361
$ws.Range("H12").Value = 11.67792797725504
$ws.Range("J12").Value = 7.282783740026254
$ws.Range("M12").Value = 20.68311830246522
$ws.Range("O12").Value = 16.42711602345656
$ws.Range("B13").Value = 15.17533875674668
$ws.Range("C13").Value = 11.31314183343794
$ws.Range("E13").Value = 23.31604585333375
$ws.Range("F13").Value = 36.47356717126299
$ws.Range("G13").Value = 18.69179728613366
$ws.Range("H13").Value = 11.68124129957446
$ws.Range("J13").Value = 7.283626813719364
$ws.Range("M13").Value = 20.66453234468177
$ws.Range("O13").Value = 16.4326138857376
$ws.Range("B14").Value = 15.05797205152887
$ws.Range("C14").Value = 11.24640044653965
$ws.Range("E14").Value = 23.30198936112272
$ws.Range("F14").Value = 36.45990939338168
$ws.Range("G14").Value = 18.70723027997059
$ws.Range("H14").Value = 11.69211216623639
$ws.Range("J14").Value = 7.286390671414159
$ws.Range("M14").Value = 20.60384409686417
$ws.Range("O14").Value = 16.45068458942214
$ws.Range("B15").Value = 14.98560553652229
$ws.Range("C15").Value = 11.20531024165319
$ws.Range("E15").Value = 23.29346375274149
$ws.Range("F15").Value = 36.45176084093762
$ws.Range("G15").Value = 18.71692561487085
$ws.Range("H15").Value = 11.69882113341272
$ws.Range("J15").Value = 7.28809469393611
$ws.Range("M15").Value = 20.5666140890073
$ws.Range("O15").Value = 16.46186164564052
$ws.Range("B16").Value = 14.56362910019799
$ws.Range("C16").Value = 10.96667376467786
$ws.Range("E16").Value = 23.24590804915148
$ws.Range("F16").Value = 36.40840031788627
$ws.Range("G16").Value = 18.77617668986612
$ws.Range("H16").Value = 11.73803730692398
$ws.Range("J16").Value = 7.298029834818817
$ws.Range("M16").Value = 20.35244639998486
$ws.Range("O16").Value = 16.52756745166797
$ws.Range("B17").Value = 14.29829764853198
$ws.Range("C17").Value = 10.81750066461241
$ws.Range("E17").Value = 23.21790309467383
$ws.Range("F17").Value = 36.38478225477412
$ws.Range("G17").Value = 18.81579571236789
$ws.Range("H17").Value = 11.76278168035321
$ws.Range("J17").Value = 7.30427675988435
$ws.Range("M17").Value = 20.22038880270481
$ws.Range("O17").Value = 16.56934771109347
$ws.Range("B18").Value = 14.14332304696403
$ws.Range("C18").Value = 10.73069920945075
$ws.Range("E18").Value = 23.20222800477938
$ws.Range("F18").Value = 36.37229577229399
$ws.Range("G18").Value = 18.83977604562689
$ws.Range("H18").Value = 11.77726613782546
$ws.Range("J18").Value = 7.307925808981734
$ws.Range("M18").Value = 20.1442084469175
$ws.Range("O18").Value = 16.59391772330025
$ws.Range("B19").Value = 14.09044671911179
$ws.Range("C19").Value = 10.70113976081696
$ws.Range("E19").Value = 23.19699527050198
$ws.Range("F19").Value = 36.36825681594525
$ws.Range("G19").Value = 18.84809929035198
$ws.Range("H19").Value = 11.7822136196368
$ws.Range("J19").Value = 7.309170939889649
$ws.Range("M19").Value = 20.11837876086607
$ws.Range("O19").Value = 16.60232913114877
$ws.Range("B20").Value = 14.32678748280897
$ws.Range("C20").Value = 10.83348445874388
$ws.Range("E20").Value = 23.2208395623408
$ws.Range("F20").Value = 36.38718284049671
$ws.Range("G20").Value = 18.81145456528044
$ws.Range("H20").Value = 11.76012149993609
$ws.Range("J20").Value = 7.303605972488522
$ws.Range("M20").Value = 20.23447029670867
$ws.Range("O20").Value = 16.56484429476257
$ws.Range("B21").Value = 15.09254475140263
$ws.Range("C21").Value = 11.26604769742267
$ws.Range("E21").Value = 23.30610052321989
$ws.Range("F21").Value = 36.46387575825413
$ws.Range("G21").Value = 18.70264671578427
$ws.Range("H21").Value = 11.68890865677411
$ws.Range("J21").Value = 7.285576550977561
$ws.Range("M21").Value = 20.62168164565027
$ws.Range("O21").Value = 16.44535422082325
$ws.Range("B22").Value = 15.57173607482108
$ws.Range("C22").Value = 11.53944144614879
$ws.Range("E22").Value = 23.36560934242288
$ws.Range("F22").Value = 36.52372771847298
$ws.Range("G22").Value = 18.64234901651709
$ws.Range("H22").Value = 11.64461550336004
$ws.Range("J22").Value = 7.27428948505064
$ws.Range("M22").Value = 20.87226976692486
$ws.Range("O22").Value = 16.37209973960148
$ws.Range("B23").Value = 15.31793892838771
$ws.Range("C23").Value = 11.39439568442562
$ws.Range("E23").Value = 23.33350534094425
$ws.Range("F23").Value = 36.49089622713831
$ws.Range("G23").Value = 18.67353232583578
$ws.Range("H23").Value = 11.6680496860161
$ws.Range("J23").Value = 7.280268313416357
$ws.Range("M23").Value = 20.73877544486436
$ws.Range("O23").Value = 16.41075236425959
$ws.Range("B24").Value = 14.31391480112283
$ws.Range("C24").Value = 10.826261419736
$ws.Range("E24").Value = 23.21951065968988
$ws.Range("F24").Value = 36.38609413411483
$ws.Range("G24").Value = 18.81341345333528
$ws.Range("H24").Value = 11.76132336203499
$ws.Range("J24").Value = 7.303909055990092
$ws.Range("M24").Value = 20.2281048535887
$ws.Range("O24").Value = 16.56687857651416
$ws.Range("B25").Value = 13.14080705524141
$ws.Range("C25").Value = 10.17548561002382
$ws.Range("E25").Value = 23.11311943779865
$ws.Range("F25").Value = 36.31503297362912
$ws.Range("G25").Value = 19.00976212231538
$ws.Range("H25").Value = 11.87157297043859
$ws.Range("J25").Value = 7.321952332463384
$ws.Range("M25").Value = 19.66899598536668
$ws.Range("O25").Value = 16.75587584400248
